# Apply "GetAllBookings im BookingRepositoryDB erstellt" update to the Tasks sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 67 ("DeleteUser Methode im AdminController erstellen") is now done:
# copy the "Zustand"/"abgeschlossen am" cells (C/D) from the row above (already
# marked done) so the new cells get the same shared string ("done") and the
# same date style as the existing entries.
$ws.Cells.Item(66, 3).Copy($ws.Cells.Item(67, 3)) | Out-Null
$ws.Cells.Item(66, 4).Copy($ws.Cells.Item(67, 4)) | Out-Null

# Row 68 used to reference the task with the typo "BookinRepositoryDB".
# It is replaced by the corrected task text, and the row is marked done too.
$ws.Cells.Item(68, 2).Value2 = "GetAllBookings im BookingRepositoryDB erstellen"
$ws.Cells.Item(66, 3).Copy($ws.Cells.Item(68, 3)) | Out-Null
$ws.Cells.Item(66, 4).Copy($ws.Cells.Item(68, 4)) | Out-Null

# Rows 69, 71 and 72 are reordered: the "(Delete bereits im BookingRepositoryDB
# vorhanden)" note now comes before the "DeleteBooking Methode" task.
$ws.Cells.Item(69, 2).Value2 = "ShowBookings Methode im AdminController erstellen"
$ws.Cells.Item(71, 2).Value2 = "(Delete bereits im BookingRepositoryDB vorhanden)"
$ws.Cells.Item(72, 2).Value2 = "DeleteBooking Methode im AdminController erstellen"

# Update the active selection to reflect the last edited cell.
$ws.Range("D68").Select() | Out-Null
